# Weekly crime-data refresh: new week header dates + volume/issue number,
# and the updated weekly/28-day/YTD/2-year crime-complaint figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: issue number and reporting week dates ---
$ws.Range("A8").Value = "Volume 31   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/29/2024  Through  2/4/2024"

# Row 15
$ws.Range("N15").Copy($ws.Range("M15"))
$ws.Range("M15").Value = -100

# Row 16
$ws.Range("D16").Copy($ws.Range("C16"))
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 0
$ws.Range("L16").Value = -80
$ws.Range("M16").Value = -90
$ws.Range("N16").Value = -96.774193548387

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 10
$ws.Range("H17").Value = 11.111111111111
$ws.Range("I17").Value = 10
$ws.Range("J17").Value = 11
$ws.Range("K17").Value = -9.090909090909
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 11.111111111111
$ws.Range("N17").Value = -28.571428571428

# Row 18
$ws.Range("D16").Copy($ws.Range("C18"))
$ws.Range("D16").Copy($ws.Range("D18"))
$ws.Range("E16").Copy($ws.Range("E18"))
$ws.Range("F18").Value = 2
$ws.Range("H18").Value = -60
$ws.Range("M18").Value = -25
$ws.Range("N18").Value = -93.617021276595

# Row 19
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 0
$ws.Range("I19").Value = 14
$ws.Range("J19").Value = 17
$ws.Range("K19").Value = -17.647058823529
$ws.Range("L19").Value = -22.222222222222
$ws.Range("M19").Value = 16.666666666666
$ws.Range("N19").Value = 40

# Row 20
$ws.Range("D16").Copy($ws.Range("C20"))
$ws.Range("D16").Copy($ws.Range("D20"))
$ws.Range("E16").Copy($ws.Range("E20"))
$ws.Range("F20").Value = 2
$ws.Range("H20").Value = 100
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = -80.95238095238

# Row 21
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 4
$ws.Range("E21").Value = 25
$ws.Range("G21").Value = 31
$ws.Range("H21").Value = -12.903225806451
$ws.Range("I21").Value = 32
$ws.Range("J21").Value = 37
$ws.Range("K21").Value = -13.513513513513
$ws.Range("L21").Value = -17.948717948717
$ws.Range("M21").Value = -13.513513513513
$ws.Range("N21").Value = -74.4

# Row 23
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = -50
$ws.Range("M23").Value = -75

# Row 24
$ws.Range("C24").Value = 4
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = -55.555555555555
$ws.Range("F24").Value = 23
$ws.Range("G24").Value = 45
$ws.Range("H24").Value = -48.888888888888
$ws.Range("I24").Value = 31
$ws.Range("J24").Value = 61
$ws.Range("K24").Value = -49.180327868852
$ws.Range("L24").Value = 0

# Row 25
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 25
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = 7.142857142857
$ws.Range("I25").Value = 20
$ws.Range("J25").Value = 16
$ws.Range("M25").Value = -13.043478260869

# Row 27
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 5
$ws.Range("K27").Value = 150
$ws.Range("L27").Value = 400

# Row 28
$ws.Range("N28").Copy($ws.Range("L28"))
$ws.Range("L28").Value = -100
$ws.Range("N28").Copy($ws.Range("M28"))
$ws.Range("M28").Value = -100

# Row 29
$ws.Range("N29").Copy($ws.Range("L29"))
$ws.Range("L29").Value = -100
$ws.Range("N29").Copy($ws.Range("M29"))
$ws.Range("M29").Value = -100
